$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.069.43"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "3.831.46"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D7").Value = "3.830.30"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "4.474.09"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.880.37"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "71.043.75"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.737"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").Value = "3.793.25"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "428.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("E51").Value = "  -2.22%  "
